$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 19 (the old aggregate/summary row with A19=17)
$ws.Rows.Item(19).Delete()

# Update B3:B18 inline string values per the relabelled dataset
$ws.Range("B3").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.563|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.437|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B4").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:1.000|x17:0.000"
$ws.Range("B5").Value = "x1:0.000|x2:0.000|x3:0.000|x4:1.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B6").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:1.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B7").Value = "x1:1.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B8").Value = "x1:0.000|x2:1.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B9").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:1.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B10").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:1.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B11").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:1.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B12").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:1.000|x16:0.000|x17:0.000"
$ws.Range("B13").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:1.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B14").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:1.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B15").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:1.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B16").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:1.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B17").Value = "x1:0.000|x2:0.000|x3:1.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B18").Value = "x1:1.000|x2:1.000|x3:2.000|x4:1.000|x5:1.000|x6:1.000|x7:0.563|x8:1.000|x9:1.000|x10:1.000|x11:1.000|x12:0.437|x13:1.000|x14:1.000|x15:1.000|x16:1.000|x17:0.000"
